$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J51").Value = 7300
$ws.Range("K51").Value = 7399
$ws.Range("L51").Value = 7300
$ws.Range("M51").Value = -6915
$ws.Range("N51").Value = -8268
$ws.Range("I51").Value = 7399
$ws.Range("K53").Value = 1893.3334
$ws.Range("H53").Value = 1479
$ws.Range("I53").Value = 1893.3334
$ws.Range("M53").Value = -1256.3334
$ws.Range("H131").Value = 1743
$ws.Range("M131").Value = -189
$ws.Range("I131").Value = 1743
$ws.Range("K131").Value = 5229
$ws.Range("H135").Value = 584.3077
$ws.Range("M135").Value = -3079.4997
$ws.Range("K135").Value = 5614.4997
$ws.Range("I135").Value = 623.8333
$ws.Range("I137").Value = 23810524
$ws.Range("K137").Value = 71431572
$ws.Range("M137").Value = -71429022
$ws.Range("H137").Value = 15874373
$ws.Range("N138").Value = -21001.5905
$ws.Range("J138").Value = 3573.8635
$ws.Range("L138").Value = 10721.5905
$ws.Range("H138").Value = 3853.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("M26").Value = $null
$ws.Range("M32").Value = -4059.9287
$ws.Range("H32").Value = 6633.3125
$ws.Range("I32").Value = 4346.9287
$ws.Range("K32").Value = 4346.9287
$ws.Range("N88").Value = $null
$ws.Range("H88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("N91").Value = $null
$ws.Range("H91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("H130").Value = 26874
$ws.Range("L130").Value = 26874
$ws.Range("J130").Value = 26874
$ws.Range("N130").Value = -36914
$ws.Range("K132").Value = 8683.2855
$ws.Range("I132").Value = 2894.4285
$ws.Range("H132").Value = 2894.4285
$ws.Range("M132").Value = -6153.2855

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M86").Value = -24211.416
$ws.Range("I86").Value = 25334.416
$ws.Range("N86").Value = -30325894
$ws.Range("K86").Value = 25334.416
$ws.Range("H86").Value = 14515833
$ws.Range("J86").Value = 30323648
$ws.Range("L86").Value = 30323648
$ws.Range("M89").Value = -121056.08
$ws.Range("I89").Value = 25334.416
$ws.Range("N89").Value = -151629472
$ws.Range("K89").Value = 126672.08
$ws.Range("H89").Value = 14515833
$ws.Range("L89").Value = 151618240
$ws.Range("J89").Value = 30323648
$ws.Range("I99").Value = 3000
$ws.Range("H99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M6").Value = -12502037
$ws.Range("I6").Value = 12502150
$ws.Range("H6").Value = 10003060
$ws.Range("K6").Value = 12502150
$ws.Range("M7").Value = -58823617
$ws.Range("I7").Value = 58823730
$ws.Range("H7").Value = 52631770
$ws.Range("L7").Value = 117.5
$ws.Range("J7").Value = 117.5
$ws.Range("N7").Value = -343.5
$ws.Range("K7").Value = 58823730
$ws.Range("M17").Value = -4015
$ws.Range("H17").Value = 5141.75
$ws.Range("J17").Value = 8000
$ws.Range("N17").Value = -8348
$ws.Range("K17").Value = 4189
$ws.Range("L17").Value = 8000
$ws.Range("I17").Value = 4189
$ws.Range("J22").Value = 490
$ws.Range("N22").Value = -1190
$ws.Range("H22").Value = 426.66666
$ws.Range("M22").Value = -45
$ws.Range("L22").Value = 490
$ws.Range("I22").Value = 395
$ws.Range("K22").Value = 395
$ws.Range("H25").Value = 5712.5
$ws.Range("M25").Value = -5776
$ws.Range("I25").Value = 5950
$ws.Range("K25").Value = 5950
$ws.Range("H31").Value = 2083.138
$ws.Range("L31").Value = 3007.6
$ws.Range("J31").Value = 3007.6
$ws.Range("N31").Value = -3597.6
$ws.Range("L34").Value = 3007.6
$ws.Range("J34").Value = 3007.6
$ws.Range("N34").Value = -3411.6
$ws.Range("H34").Value = 2083.138
$ws.Range("J41").Value = 30000
$ws.Range("H41").Value = 28511.8
$ws.Range("N41").Value = -30856
$ws.Range("L41").Value = 30000
$ws.Range("M62").Value = -1776
$ws.Range("N62").Value = -3747.5
$ws.Range("I62").Value = 2400
$ws.Range("H62").Value = 2439.8
$ws.Range("K62").Value = 2400
$ws.Range("J62").Value = 2499.5
$ws.Range("L62").Value = 2499.5
$ws.Range("L65").Value = 12497.5
$ws.Range("N65").Value = -18737.5
$ws.Range("J65").Value = 2499.5
$ws.Range("M65").Value = -8880
$ws.Range("I65").Value = 2400
$ws.Range("H65").Value = 2439.8
$ws.Range("K65").Value = 12000
$ws.Range("H93").Value = 13675.4
$ws.Range("I93").Value = 13675.4
$ws.Range("M93").Value = -11803.4
$ws.Range("K93").Value = 13675.4
$ws.Range("I99").Value = 4736
$ws.Range("H99").Value = 13108.833
$ws.Range("K99").Value = 4736
$ws.Range("M99").Value = -3238
$ws.Range("M126").Value = -11738
$ws.Range("K126").Value = 14208
$ws.Range("I126").Value = 4736
$ws.Range("H126").Value = 13108.833
$ws.Range("K132").Value = 27292947
$ws.Range("I132").Value = 9097649
$ws.Range("H132").Value = 8702177
$ws.Range("M132").Value = -27290417
$ws.Range("H134").Value = 4097.8
$ws.Range("I134").Value = 4097.8
$ws.Range("K134").Value = 12293.4
$ws.Range("M134").Value = -9758.400000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M2").Value = -213.4
$ws.Range("H2").Value = 571464.6
$ws.Range("K2").Value = 326.4
$ws.Range("I2").Value = 54.4
$ws.Range("M7").Value = -390366.008
$ws.Range("I7").Value = 130159.336
$ws.Range("H7").Value = 97683.836
$ws.Range("L7").Value = 772.0000200000001
$ws.Range("J7").Value = 257.33334
$ws.Range("N7").Value = -996.0000200000001
$ws.Range("K7").Value = 390478.008
$ws.Range("M11").Value = -4725.25
$ws.Range("K11").Value = 4865.25
$ws.Range("I11").Value = 1621.75
$ws.Range("H11").Value = 1431.2
$ws.Range("N39").Value = -23726.1432
$ws.Range("L39").Value = 23138.1432
$ws.Range("J39").Value = 7712.7144
$ws.Range("H39").Value = 7811.125
$ws.Range("M46").Value = -2009
$ws.Range("J46").Value = 957.8
$ws.Range("L46").Value = 2873.4
$ws.Range("H46").Value = 914.8333
$ws.Range("I46").Value = 700
$ws.Range("K46").Value = 2100
$ws.Range("N46").Value = -3055.4
$ws.Range("N117").Value = -10067
$ws.Range("H117").Value = 882.6667
$ws.Range("J117").Value = 1061
$ws.Range("L117").Value = 3183
$ws.Range("H134").Value = 5325.2856
$ws.Range("I134").Value = 5325.2856
$ws.Range("K134").Value = 15975.8568
$ws.Range("M134").Value = -10905.8568

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I97").Value = 1570.2
$ws.Range("H97").Value = 1534.5625
$ws.Range("K97").Value = 1570.2
$ws.Range("M97").Value = -1074.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3498.5
$ws.Range("L7").Value = 2999
$ws.Range("J7").Value = 2999
$ws.Range("N7").Value = -3223
$ws.Range("J22").Value = 500000000
$ws.Range("N22").Value = -500000590
$ws.Range("H22").Value = 500000000
$ws.Range("M22").Value = $null
$ws.Range("L22").Value = 500000000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = $null
$ws.Range("H27").Value = 500000000
$ws.Range("J27").Value = 500000000
$ws.Range("N27").Value = -500000214
$ws.Range("L27").Value = 500000000
$ws.Range("K61").Value = 3157.8667
$ws.Range("H61").Value = 3112.9375
$ws.Range("M61").Value = -2955.8667
$ws.Range("I61").Value = 3157.8667
$ws.Range("M82").Value = -508
$ws.Range("K82").Value = 869
$ws.Range("I82").Value = 869
$ws.Range("H82").Value = 1229.8462
$ws.Range("K85").Value = 869
$ws.Range("I85").Value = 869
$ws.Range("H85").Value = 1229.8462
$ws.Range("M85").Value = 379
$ws.Range("I100").Value = 1699
$ws.Range("M100").Value = -1158
$ws.Range("H100").Value = 1699
$ws.Range("K100").Value = 1699
$ws.Range("K113").Value = 3157.8667
$ws.Range("M113").Value = -987.8667
$ws.Range("I113").Value = 3157.8667
$ws.Range("H113").Value = 3112.9375
$ws.Range("L126").Value = 8997
$ws.Range("J126").Value = 2999
$ws.Range("N126").Value = -13937
$ws.Range("H126").Value = 3498.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J81").Value = 10042.889
$ws.Range("N81").Value = -22207.778
$ws.Range("I81").Value = 3110.2
$ws.Range("K81").Value = 6220.4
$ws.Range("L81").Value = 20085.778
$ws.Range("H81").Value = 5261.724
$ws.Range("M81").Value = -5159.4
$ws.Range("L84").Value = 100428.89
$ws.Range("H84").Value = 5261.724
$ws.Range("M84").Value = -25798
$ws.Range("J84").Value = 10042.889
$ws.Range("I84").Value = 3110.2
$ws.Range("K84").Value = 31102
$ws.Range("N84").Value = -111036.89
$ws.Range("L132").Value = 750004500
$ws.Range("K132").Value = 4900.8462
$ws.Range("J132").Value = 250001500
$ws.Range("I132").Value = 1633.6154
$ws.Range("H132").Value = 58825132
$ws.Range("M132").Value = -2370.8462
$ws.Range("N132").Value = -750009560
